# Weekly data refresh: a new week's record is inserted at the top of the
# data block (row 587), pushing all the existing records (previously rows
# 587-632) down by one row (to 588-633). This grows the used range from
# A1:T632 to A1:T633.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 587; Excel shifts 587:632 down to 588:633
# and extends the sheet dimension automatically.
$ws.Rows("587:587").Insert()

# Populate the newly inserted row with the new week's record.
$ws.Range("A587").Value = 9
$ws.Range("B587").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C587").Value = "Metropolitana"
$ws.Range("D587").Value = 45013
$ws.Range("E587").Value = 13
$ws.Range("F587").Value = "Fruta"
$ws.Range("G587").Value = 100108
$ws.Range("H587").Value = "Tropicales y subtropicales"
$ws.Range("I587").Value = 100108002
$ws.Range("J587").Value = "Mango"
$ws.Range("K587").Value = "Sin especificar"
$ws.Range("L587").Value = "Primera"
$ws.Range("M587").Value = 730
$ws.Range("N587").Value = 6500
$ws.Range("O587").Value = 7000
$ws.Range("P587").Value = 6801
$ws.Range("Q587").Value = "$/bandeja 4 kilos"
$ws.Range("R587").Value = "Perú"
$ws.Range("S587").Value = 1700
$ws.Range("T587").Value = 4
